$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.28244960308075
$ws.Range("B1").Value = 4.815610408782959
$ws.Range("C1").Value = 3.266622066497803
$ws.Range("D1").Value = 1.722324967384338
$ws.Range("E1").Value = 1.290124416351318
